$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 3 and 4 entirely (they contained extra Q&A pairs no longer needed)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Update the remaining question/answer pair in row 2
$ws.Range("B2").Value = "الاستاذ محمد"
$ws.Range("A2").Value = "من رئيس جريدة يونا "

# Set explicit column widths (values chosen so the resulting stored width
# matches the target as closely as this runtime's width quantization allows)
$ws.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17.5

# Update selection to A2
$ws.Range("A2").Select()
